# Apply crypto price/volume updates per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.256.65"
$ws.Range("E2").Value = "  +2.53%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.060.35"
$ws.Range("E3").Value = "  +1.37%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.54"
$ws.Range("E5").Value = "  -0.30%  "

$ws.Range("E6").Value = "  +0.59%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.85"
$ws.Range("E7").Value = "  +8.10%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  +1.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0830"
$ws.Range("E10").Value = "  +6.01%  "

$ws.Range("E11").Value = "  +1.45%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.83"
$ws.Range("E12").Value = "  +3.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.360.61"
$ws.Range("E13").Value = "  +1.46%  "

$ws.Range("E14").Value = "  +3.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.763"
$ws.Range("E15").Value = "  +2.73%  "

$ws.Range("E16").Value = "  +2.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.044.96"
$ws.Range("E17").Value = "  +0.77%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.152.98"
$ws.Range("E18").Value = "  +2.48%  "

$ws.Range("E19").Value = "  +0.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.81"
$ws.Range("E20").Value = "  +1.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0833"
$ws.Range("E21").Value = "  +1.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.28"
$ws.Range("E22").Value = "  +0.80%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("E24").Value = "  -0.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.21"
$ws.Range("E25").Value = "  -1.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.52"
$ws.Range("E26").Value = "  +1.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.28"
$ws.Range("E27").Value = "  +1.81%  "

$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.00"
$ws.Range("E29").Value = "  +1.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.30"
$ws.Range("E30").Value = "  -1.70%  "

$ws.Range("E31").Value = "  +2.61%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.49"
$ws.Range("E32").Value = "  +0.47%  "

$ws.Range("E33").Value = "  +2.60%  "

$ws.Range("E34").Value = "  +2.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0607"
$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.28"
$ws.Range("E36").Value = "  +12.76%  "

$ws.Range("E37").Value = "  -2.18%  "

$ws.Range("E38").Value = "  +1.02%  "

$ws.Range("E39").Value = "  +0.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.528.94"
$ws.Range("E40").Value = "  +3.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.10"
$ws.Range("E41").Value = "  +3.93%  "

$ws.Range("E42").Value = "  +1.42%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.84"
$ws.Range("E43").Value = "  +3.37%  "

$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0931"
$ws.Range("E44").Value = "  +1.07%  "

$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.83"
$ws.Range("E45").Value = "  +0.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.12"
$ws.Range("E46").Value = "  +1.41%  "

$ws.Range("E47").Value = "  -6.84%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.02"
$ws.Range("E48").Value = "  +0.36%  "

$ws.Range("E49").Value = "  +1.77%  "

$ws.Range("E50").Value = "  -0.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.251.54"
$ws.Range("E51").Value = "  +1.73%  "
